$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9616913199424744
$ws.Range("B1").Value = 2.712293148040771
$ws.Range("C1").Value = 4.768298149108887
$ws.Range("D1").Value = 1.251917600631714
$ws.Range("E1").Value = 1.304546475410461
